$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 435, shifting existing rows 435-520 down to 436-521.
$ws.Rows.Item(435).Insert()

# Populate the new row 435 with the new weekly data record.
$ws.Cells.Item(435, 1).Value = 3
$ws.Cells.Item(435, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(435, 3).Value = "Coquimbo"
$ws.Cells.Item(435, 4).Value = 45015
$ws.Cells.Item(435, 5).Value = 5
$ws.Cells.Item(435, 6).Value = 100112009
$ws.Cells.Item(435, 7).Value = "Acelga"
$ws.Cells.Item(435, 8).Value = "Sin especificar"
$ws.Cells.Item(435, 9).Value = "Primera"
$ws.Cells.Item(435, 10).Value = 210
$ws.Cells.Item(435, 11).Value = 3500
$ws.Cells.Item(435, 12).Value = 3800
$ws.Cells.Item(435, 13).Value = 3643
$ws.Cells.Item(435, 14).Value = "`$/docena de atados (6 kilos)"
$ws.Cells.Item(435, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(435, 16).Value = 607
$ws.Cells.Item(435, 17).Value = 6
$ws.Cells.Item(435, 18).Value = "Hortaliza"
